$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: numeric index 0..3 (was shared-string "Name"/"Phone" before, now plain
# numbers spanning four rows)
$ws.Range("A1").Value = 0
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Column B: the old "Name"/"Phone" labels move over to column B
$ws.Range("B1").Value = "Name"
$ws.Range("B2").Value = "Phone"

# Column C: the old "Mamun"/"01643091606" values move over to column C.
# The phone number must stay text (leading zero) like the original shared
# string, so force a Text format before typing it, then drop back to the
# workbook's default (Normal) style so no stray per-cell formatting is left
# behind - only the cell's string type needs to persist.
$ws.Range("C1").Value = "Mamun"
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "01643091606"
$ws.Cells.Item(2, 3).Style = "Normal"

# Match the saved selection/page setup seen in the target file
$ws.PageSetup.Orientation = 1
$ws.Range("C5").Select()
